# Build site at 2022-09-26 16:07:08 UTC
# Rework the LOQ4059 syllabus sheet: rows 10-24 get new content/ordering,
# and the sheet shrinks by one row (old row 24 disappears).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out the old rows 10-24 completely (no leftover empty <c> placeholders)
# so the rebuilt rows below end up with exactly the cells the target has.
$ws.Rows("10:24").Delete()

# --- row 10 --------------------------------------------------------------
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "5840772 - Amilton Martins dos Santos"
$ws.Range("C10").Value = "5840772 - Amilton Martins dos Santos"
$ws.Rows.Item(10).RowHeight = 60

# --- row 11 ----------------------------------------------------------------
$ws.Range("A11").Value = "Objectives:"
$ws.Rows.Item(11).RowHeight = 60

# --- row 12 ----------------------------------------------------------------
$ws.Range("A12").Value = "Docentes responsáveis:"

# --- row 13 ----------------------------------------------------------------
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# --- row 14 ----------------------------------------------------------------
$ws.Range("A14").Value = "Short syllabus:"
$ws.Rows.Item(14).RowHeight = 60

# --- row 15 ----------------------------------------------------------------
# "01/01/2012" looks like a date, so force text format first or Excel will
# silently convert it to a date serial number.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "01/01/2012"
$ws.Rows.Item(15).RowHeight = 120

# --- row 16 ----------------------------------------------------------------
$ws.Range("A16").Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

# --- row 17 ----------------------------------------------------------------
$ws.Range("A17").Value = "Avaliação:"

# --- row 18 ----------------------------------------------------------------
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840772 - Amilton Martins dos Santos"
$ws.Range("C18").Value = "5840772 - Amilton Martins dos Santos"
$ws.Rows.Item(18).RowHeight = 60

# --- row 19 ----------------------------------------------------------------
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "2 Provas escritas + Trabalho de conclusão de curso."
$ws.Range("C19").Value = "2 Provas escritas + Trabalho de conclusão de curso."
$ws.Rows.Item(19).RowHeight = 60

# --- row 20 ----------------------------------------------------------------
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota final (NF) será calculada de seguinte maneira: NF = (P1+P2)/2`nO trabalho poderá valer até 2 pontos, que serão somados nas notas da P1 ou da P2."
$ws.Range("C20").Value = "A nota final (NF) será calculada de seguinte maneira: NF = (P1+P2)/2`nO trabalho poderá valer até 2 pontos, que serão somados nas notas da P1 ou da P2."
$ws.Rows.Item(20).RowHeight = 60

# --- row 21 ----------------------------------------------------------------
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula MR = (NF+PR)/2."
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula MR = (NF+PR)/2."
$ws.Rows.Item(21).RowHeight = 120

# --- row 22 ----------------------------------------------------------------
$ws.Range("A22").Value = "Requisitos:"

# --- row 23 ----------------------------------------------------------------
$ws.Range("B23").Value = "LOQ4038 -  Química Orgânica II  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOQ4038 -  Química Orgânica II  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30

# --- fix up column-B styling -----------------------------------------------
# Freshly-created column-B cells pick up the bold "column A" style (1)
# instead of the wrap-text style (2) the B column actually uses, because the
# <cols> definition has an overlapping min=1/max=2 range ahead of the
# min=2/max=2 one. Re-apply the correct format from an untouched B cell.
$ws.Range("B3").Copy()
foreach ($r in 10,13,15,18,19,20,21,23) {
    $ws.Range("B$r").PasteSpecial(-4122)
}

# Row 15's NumberFormat="@" override (needed so "01/01/2012" stays text)
# also leaves column C with a stray one-off style; re-sync it too.
$ws.Range("C3").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$excel.CutCopyMode = $false
